# Convert the "tickersToAddTable" workbook to an "attended" version for the
# UiPath Marketplace listing:
#   - The Type/Ticker sample row (stock / GME) is removed from Sheet1, only
#     the header row remains, and column B is resized to fit its new content.
#   - The hidden lookup values on the "Data" sheet are re-cased to
#     "Stock"/"Crypto" (title case) instead of "stock"/"crypto".
#   - The "Data" sheet (used only to back the dropdown data validation on
#     Sheet1) is hidden from end users.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Data")

# --- Data sheet: re-case the lookup/validation source values -----------
$ws2.Range("A1").Value2 = "Stock"
$ws2.Range("A2").Value2 = "Crypto"

# --- Sheet1: drop the example row (stock / GME), keep only the header --
$ws1.Rows.Item(2).Delete()

# Resize column B to fit the (now shorter) header-only content
$ws1.Columns.Item(2).AutoFit()

# Leave the cursor on the next free entry row, as the author did
[void]$ws1.Range("B3").Select()

# --- Hide the helper "Data" sheet so only Sheet1 is shown to users -----
$ws2.Visible = $false
